$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13-19 down to 14-20)
$ws.Rows.Item(13).Insert()

# Copy formatting from row 12 into new row 13
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(13).PasteSpecial(-4122)  # xlPasteFormats

# Set values for the new row
$ws.Range("A13").Value = "D6"
$ws.Range("B13").Value = "LevelSense"
$ws.Range("C13").Value = "Generic (D6)"
$ws.Range("D13").Value = "Levelsense"

$ws.Range("D14").Select()
